$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$messages = @(
    "Ingen togdrift i spor 1. Passagerer til Orientkaj, skal skifte tog på Østerport.",
    "M3/M4: Alle tog holder stille grunde tekniske problemer. Rejsegarantien er trådt i kraft.",
    "Vi holder stille i begge spor grundet tekniske problemer.",
    "Ingen tog standser i øjeblikket på Poul Henningsens Plads.",
    "Ingen tog standser i spor 1, på Amager Strand, retning CPH Lufthavnen. Skal du af på Amager Strand, kan du køre til Femøren og tage toget tilbage til Amager Strand. Skal du fra Amager Strand, mod CPH Lufthavnen kan du tage til Øresund og tage toget mod CPH Lufthavnen. Lidt længere rejsetid må forventes",
    "Længere rejsetid og togskift på flere stationer i spor 1 på M3 linjen. Der er normal drift i spor 2. Der er togskift på Østerport og Nordhavn på M4 linjen. Det skyldes tekniske problemer med et sporskifte.",
    "M.dk sat i for M3/M4: Passagerer der kører i spor 1 skal skifte tog på Østerport/ Gl. Strand/Enghave Plads/Axel Møllers Have, Nørrebro og Trianglen. Det skyldes problemer med et sporskifte. Forvent derfor længere rejsetid.`nVi kører normalt i spor 2",
    "M2: Ingen togdrift mellem Christianshavn - Lufthavnen. Skyldes et tog med en teknisk fejl ved Lergravsparken. Rejsegarantien er aktiv mellem Christianshavn - Lufthavnen.",
    "M2: Forsinkelser særligt mellem Christianshavn - Lufthavnen pga. et tog med en teknisk fejl ved Lergravsparken.",
    "Ingen togdrift mellem Vanløse - Frederiksberg. Skyldes politiets arbejde nær Lindevang.",
    "Ingen togdrift mellem Vanløse - Frederiksberg. Dette skyldes politiets arbejde nær Lindevang",
    "Ingen togdrift mellem Vanløse - Frederiksberg. Dette skyldes politiets arbejde"
)

for ($i = 0; $i -lt $messages.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $messages[$i]
}
